# Applies the "update SSR and TR" commit to the TR rules table on Sheet1.
#
# Summary of the change:
#   - B11 ("OPassVS") is renamed to "SPassPartVByS" and its look is changed
#     from the Arial/10 style to the Helvetica/10 style (matches the other
#     rows further down the table).
#   - Seven new transformation rules (TR26..TR32) are appended as rows 27-33,
#     reusing the existing "Helvetica/10" (rows mostly) and "Arial/10"
#     (the add_relation / add_entity_attribute rows) cell styles already
#     present in the workbook, except for the very last new cell (B33) which
#     introduces a brand-new font (Calibri/11, black, minor scheme).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Style donor cells already present in the workbook - copying *formats*
# from them (instead of touching Font.* directly) guarantees the new
# cells land on the very same cellXfs/font entries instead of Excel
# minting brand-new (duplicate) style records.
# ---------------------------------------------------------------------
$styleArial10     = $ws.Range("C5")   # s="2"  Arial 10
$styleCalibri11   = $ws.Range("B2")   # s="3"  Calibri 11 black
$styleHelvetica10 = $ws.Range("B19")  # s="4"  Helvetica 10

function Set-CellFormat($cell, $donor) {
    $donor.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------
# 1) B11: OPassVS -> SPassPartVByS, style Arial10 -> Helvetica10
# ---------------------------------------------------------------------
$b11 = $ws.Range("B11")
Set-CellFormat $b11 $styleHelvetica10
$b11.Value = "SPassPartVByS"

# ---------------------------------------------------------------------
# 2) New rows 27-33 (TR26 .. TR32)
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=27; A="TR26"; B="SVIOO";          Bs="H"; C="add_behavior";      Cs="H"; D='actor=B, target=D, action=A' }
    @{ Row=28; A="TR27"; B="SVIOO";          Bs="H"; C="add_relation";      Cs="A"; D='source=D, dest=C, msg="to", ass_type="association"' }
    @{ Row=29; A="TR28"; B="SVOTo";          Bs="H"; C="add_behavior";      Cs="H"; D='actor=B, target=C, action=A' }
    @{ Row=30; A="TR29"; B="SVOTo";          Bs="H"; C="add_relation";      Cs="A"; D='source=C, dest=D, msg="to", ass_type="association"' }
    @{ Row=31; A="TR30"; B="OPassPartVByS";  Bs="H"; C="add_behavior";      Cs="H"; D='actor=C, target=B, action=A' }
    @{ Row=32; A="TR31"; B="SVOPassPart";    Bs="A"; C="add_behavior";      Cs="A"; D='actor=B, target=D, action=A' }
    @{ Row=33; A="TR32"; B="SPredicativeV";  Bs="N"; C="add_relation";      Cs="H"; D='source=A, dest=B, ass_type="generalization"' }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)
    $cellD = $ws.Cells.Item($r, 4)

    # --- column B style ---
    if ($item.Bs -eq "H") {
        Set-CellFormat $cellB $styleHelvetica10
    } elseif ($item.Bs -eq "A") {
        Set-CellFormat $cellB $styleArial10
    } elseif ($item.Bs -eq "N") {
        # brand-new font: Calibri 11, black, minor scheme (closest reachable
        # match via the exposed Font object model)
        $cellB.Font.Name = "Calibri"
        $cellB.Font.Size = 11
        $cellB.Font.ColorIndex = 1
    }

    # --- column C style ---
    if ($item.Cs -eq "H") {
        Set-CellFormat $cellC $styleHelvetica10
    } elseif ($item.Cs -eq "A") {
        Set-CellFormat $cellC $styleArial10
    }

    # --- values (columns A and D keep the default/no explicit style) ---
    $cellA.Value = $item.A
    $cellB.Value = $item.B
    $cellC.Value = $item.C
    $cellD.Value = $item.D
}

# ---------------------------------------------------------------------
# 3) Restore the saved cursor/selection position
# ---------------------------------------------------------------------
$ws.Range("C9").Select() | Out-Null

Write-Output "TR table updated: B11 renamed, TR26-TR32 appended."
